# This script applies a row-level permutation to rows 2-41 of the sheet,
# affecting columns D, L, M, N, O, P, Q, R, S, T (Fecha, Calidad, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg / unidad).
#
# Mapping below: for each destination row (2-41), it lists which row's
# original (pre-edit) data should end up there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","L","M","N","O","P","Q","R","S","T")

$mapping = @{
    2  = 40
    3  = 20
    4  = 21
    5  = 29
    6  = 30
    7  = 31
    8  = 23
    9  = 17
    10 = 27
    11 = 16
    12 = 34
    13 = 19
    14 = 32
    15 = 41
    16 = 9
    17 = 7
    18 = 3
    19 = 25
    20 = 26
    21 = 2
    22 = 4
    23 = 5
    24 = 24
    25 = 33
    26 = 10
    27 = 37
    28 = 15
    29 = 12
    30 = 13
    31 = 28
    32 = 14
    33 = 11
    34 = 38
    35 = 39
    36 = 6
    37 = 8
    38 = 18
    39 = 35
    40 = 22
    41 = 36
}

# Snapshot original values for every relevant cell before any writes,
# so that overlapping reads/writes during the permutation don't clobber
# data that is still needed for a later row.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    foreach ($c in $cols) {
        $ref = "$c$r"
        $snapshot[$ref] = $ws.Range($ref).Value2
    }
}

# Write back the permuted values.
for ($r = 2; $r -le 41; $r++) {
    $srcRow = $mapping[$r]
    if ($srcRow -eq $r) {
        continue
    }
    foreach ($c in $cols) {
        $srcRef = "$c$srcRow"
        $dstRef = "$c$r"
        $ws.Range($dstRef).Value = $snapshot[$srcRef]
    }
}
